$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with P1/Q1, matching style of existing header cells (O1)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Row 2
$ws.Range("B2").Value = 24.3198244938717
$ws.Range("C2").Value = 19.92594805995155
$ws.Range("D2").Value = 9.354665912739412
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 28.17812535241687
$ws.Range("G2").Value = 32.00001100594942
$ws.Range("H2").Value = 1.850882370739265
$ws.Range("I2").Value = 2.77322340408848
$ws.Range("J2").Value = 10.60496633694934
$ws.Range("K2").Value = 14.25287775458085
$ws.Range("L2").Value = 9.487458098544385
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0

# Row 3
$ws.Range("B3").Value = 22.7055870629609
$ws.Range("C3").Value = 18.61090311934585
$ws.Range("D3").Value = 8.858034749070674
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 27.21656456936723
$ws.Range("G3").Value = 30.80598228312593
$ws.Range("H3").Value = 1.564788124491925
$ws.Range("I3").Value = 2.559047919363895
$ws.Range("J3").Value = 10.49893516042273
$ws.Range("K3").Value = 14.35941808367057
$ws.Range("L3").Value = 9.076612974205034
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0

# Row 4
$ws.Range("B4").Value = 21.65464224244831
$ws.Range("C4").Value = 17.76151699475226
$ws.Range("D4").Value = 8.540661385739256
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 26.61871277579857
$ws.Range("G4").Value = 30.06698031436053
$ws.Range("H4").Value = 1.742315374356472
$ws.Range("I4").Value = 2.521081782043886
$ws.Range("J4").Value = 10.43777263327796
$ws.Range("K4").Value = 14.43157967250953
$ws.Range("L4").Value = 8.81455137159113
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0

# Row 5
$ws.Range("B5").Value = 21.21108600060113
$ws.Range("C5").Value = 17.41736900793596
$ws.Range("D5").Value = 8.412605359327822
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 26.35673957890528
$ws.Range("G5").Value = 29.73418246562628
$ws.Range("H5").Value = 1.819955321109024
$ws.Range("I5").Value = 2.582870952699339
$ws.Range("J5").Value = 10.40812357358524
$ws.Range("K5").Value = 14.45176234622248
$ws.Range("L5").Value = 8.705179627096738
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0

# Row 6
$ws.Range("B6").Value = 21.13630413848963
$ws.Range("C6").Value = 17.37495485858654
$ws.Range("D6").Value = 8.396464407622107
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 26.29267145002883
$ws.Range("G6").Value = 29.64152707151721
$ws.Range("H6").Value = 1.833248102471329
$ws.Range("I6").Value = 2.597325563823794
$ws.Range("J6").Value = 10.39630383647976
$ws.Range("K6").Value = 14.44201453024768
$ws.Range("L6").Value = 8.686663066643394
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 0

# Row 7
$ws.Range("B7").Value = 21.64814511361456
$ws.Range("C7").Value = 17.79787602100199
$ws.Range("D7").Value = 8.553153973804342
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 26.55991464686543
$ws.Range("G7").Value = 29.96179920459815
$ws.Range("H7").Value = 1.744177406717946
$ws.Range("I7").Value = 2.532896016015416
$ws.Range("J7").Value = 10.41845162602612
$ws.Range("K7").Value = 14.39576184470376
$ws.Range("L7").Value = 8.812520336286306
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 0

# Row 8
$ws.Range("B8").Value = 23.77487995150522
$ws.Range("C8").Value = 19.53114352335395
$ws.Range("D8").Value = 9.203642350690776
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 27.77823764952491
$ws.Range("G8").Value = 31.4625393284321
$ws.Range("H8").Value = 1.752443536819445
$ws.Range("I8").Value = 2.702624256558651
$ws.Range("J8").Value = 10.54293249853239
$ws.Range("K8").Value = 14.24026402715146
$ws.Range("L8").Value = 9.347175047588935
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 0

# Row 9
$ws.Range("B9").Value = 27.47545151973592
$ws.Range("C9").Value = 22.53444392050763
$ws.Range("D9").Value = 10.35522037445821
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 30.19089389271461
$ws.Range("G9").Value = 34.51445060182674
$ws.Range("H9").Value = 2.435935754680106
$ws.Range("I9").Value = 3.218293591118375
$ws.Range("J9").Value = 10.85273995624973
$ws.Range("K9").Value = 14.05935189760684
$ws.Range("L9").Value = 10.31674668831404
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 0

# Row 10
$ws.Range("B10").Value = 29.90855429315614
$ws.Range("C10").Value = 24.51300744907662
$ws.Range("D10").Value = 11.0212972190077
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 31.57008232552718
$ws.Range("G10").Value = 36.28401882498481
$ws.Range("H10").Value = 2.890404732430924
$ws.Range("I10").Value = 3.578370695450363
$ws.Range("J10").Value = 11.01605813055369
$ws.Range("K10").Value = 13.83437061511037
$ws.Range("L10").Value = 10.79907788471704
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = 0
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 0

# Row 11
$ws.Range("B11").Value = 30.92551086947309
$ws.Range("C11").Value = 25.02694644909735
$ws.Range("D11").Value = 10.28731574717342
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 29.40151692660568
$ws.Range("G11").Value = 33.69725060512669
$ws.Range("H11").Value = 3.544080294368799
$ws.Range("I11").Value = 3.683274504794857
$ws.Range("J11").Value = 10.4037134908826
$ws.Range("K11").Value = 12.68206749900259
$ws.Range("L11").Value = 9.524342343108119
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 0

# Row 12
$ws.Range("B12").Value = 31.29140472465317
$ws.Range("C12").Value = 25.00790548950012
$ws.Range("D12").Value = 9.514510013934945
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 27.32915660722096
$ws.Range("G12").Value = 31.23457836320481
$ws.Range("H12").Value = 4.605719599680926
$ws.Range("I12").Value = 3.698331318132976
$ws.Range("J12").Value = 9.87559123413724
$ws.Range("K12").Value = 11.842797294319
$ws.Range("L12").Value = 8.437592608107662
$ws.Range("M12").Value = 0
$ws.Range("N12").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0

# Row 13
$ws.Range("B13").Value = 31.19069171328421
$ws.Range("C13").Value = 24.62537609895638
$ws.Range("D13").Value = 8.663947664810594
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 25.09250054832687
$ws.Range("G13").Value = 28.52123394231891
$ws.Range("H13").Value = 5.79825338271008
$ws.Range("I13").Value = 3.64801557126056
$ws.Range("J13").Value = 9.350201832613827
$ws.Range("K13").Value = 11.15399585342806
$ws.Range("L13").Value = 7.449513635698857
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0

# Row 14
$ws.Range("B14").Value = 30.90812216211334
$ws.Range("C14").Value = 24.19763903210919
$ws.Range("D14").Value = 8.030886095506279
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 23.42787472303157
$ws.Range("G14").Value = 26.45877984548403
$ws.Range("H14").Value = 6.681337994546441
$ws.Range("I14").Value = 3.584540731739207
$ws.Range("J14").Value = 8.981634132274573
$ws.Range("K14").Value = 10.75356163482378
$ws.Range("L14").Value = 6.848601735033219
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0

# Row 15
$ws.Range("B15").Value = 30.73809201490075
$ws.Range("C15").Value = 24.02638925781752
$ws.Range("D15").Value = 7.861847777024117
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 22.97141449522108
$ws.Range("G15").Value = 25.87022453773545
$ws.Range("H15").Value = 6.881639667354627
$ws.Range("I15").Value = 3.555597769965389
$ws.Range("J15").Value = 8.887579596238037
$ws.Range("K15").Value = 10.67745892490488
$ws.Range("L15").Value = 6.712636765225959
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0

# Row 16
$ws.Range("B16").Value = 29.76206383022405
$ws.Range("C16").Value = 23.28738044082352
$ws.Range("D16").Value = 7.707004497259172
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 22.71324202533472
$ws.Range("G16").Value = 25.45072731740084
$ws.Range("H16").Value = 6.593097855067747
$ws.Range("I16").Value = 3.414446108573872
$ws.Range("J16").Value = 8.908304235850368
$ws.Range("K16").Value = 10.91862258041175
$ws.Range("L16").Value = 6.660236811362289
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0

# Row 17
$ws.Range("B17").Value = 29.15352127433501
$ws.Range("C17").Value = 22.94877532427929
$ws.Range("D17").Value = 7.936369172838755
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 23.42045116338661
$ws.Range("G17").Value = 26.26558755979865
$ws.Range("H17").Value = 5.848814913914791
$ws.Range("I17").Value = 3.33961826274519
$ws.Range("J17").Value = 9.119335051544391
$ws.Range("K17").Value = 11.29920686029128
$ws.Range("L17").Value = 6.921733717359291
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0

# Row 18
$ws.Range("B18").Value = 28.80722850727889
$ws.Range("C18").Value = 22.89959771372019
$ws.Range("D18").Value = 8.525407022285187
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 25.07953804643218
$ws.Range("G18").Value = 28.29709151277346
$ws.Range("H18").Value = 4.677244390734325
$ws.Range("I18").Value = 3.312961567752629
$ws.Range("J18").Value = 9.533474884414039
$ws.Range("K18").Value = 11.89815676835435
$ws.Range("L18").Value = 7.599682200856102
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 0

# Row 19
$ws.Range("B19").Value = 28.70322609019489
$ws.Range("C19").Value = 23.14941441561623
$ws.Range("D19").Value = 9.384807388232101
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 27.29098201749662
$ws.Range("G19").Value = 30.97309159301954
$ws.Range("H19").Value = 3.437241988208421
$ws.Range("I19").Value = 3.337342464641049
$ws.Range("J19").Value = 10.05602662206595
$ws.Range("K19").Value = 12.62667720598505
$ws.Range("L19").Value = 8.680127346387204
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = 0
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 0

# Row 20
$ws.Range("B20").Value = 29.28678508045795
$ws.Range("C20").Value = 24.10225420623086
$ws.Range("D20").Value = 10.8797100216801
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 31.04946408926318
$ws.Range("G20").Value = 35.54040434167411
$ws.Range("H20").Value = 2.768481540810103
$ws.Range("I20").Value = 3.491313001123333
$ws.Range("J20").Value = 10.91307665541123
$ws.Range("K20").Value = 13.77289626585588
$ws.Range("L20").Value = 10.66525405644004
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 0
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 0

# Row 21
$ws.Range("B21").Value = 31.07197835310654
$ws.Range("C21").Value = 25.60396501466092
$ws.Range("D21").Value = 11.55152037945971
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 32.61533810283243
$ws.Range("G21").Value = 37.53572957021857
$ws.Range("H21").Value = 3.145423931495098
$ws.Range("I21").Value = 3.773574486114644
$ws.Range("J21").Value = 11.17786523096309
$ws.Range("K21").Value = 13.82424156319085
$ws.Range("L21").Value = 11.29832417177811
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 0

# Row 22
$ws.Range("B22").Value = 32.18782070978823
$ws.Range("C22").Value = 26.49021120320903
$ws.Range("D22").Value = 11.90213163463424
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 33.5209088509779
$ws.Range("G22").Value = 38.73332261527258
$ws.Range("H22").Value = 3.376779013112608
$ws.Range("I22").Value = 3.950163501910136
$ws.Range("J22").Value = 11.33805448732932
$ws.Range("K22").Value = 13.85538732437754
$ws.Range("L22").Value = 11.61070387478279
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 0

# Row 23
$ws.Range("B23").Value = 31.59734631781177
$ws.Range("C23").Value = 25.98560992039801
$ws.Range("D23").Value = 11.70337317029179
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 33.09172235717351
$ws.Range("G23").Value = 38.18982745023514
$ws.Range("H23").Value = 3.25439936599834
$ws.Range("I23").Value = 3.853529533074961
$ws.Range("J23").Value = 11.27252549245237
$ws.Range("K23").Value = 13.88028564894235
$ws.Range("L23").Value = 11.44532920622265
$ws.Range("M23").Value = 0
$ws.Range("N23").Value = 0
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 0

# Row 24
$ws.Range("B24").Value = 29.26042313842224
$ws.Range("C24").Value = 24.05717056664361
$ws.Range("D24").Value = 10.9462140677902
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 31.3612057095256
$ws.Range("G24").Value = 35.96387692774169
$ws.Range("H24").Value = 2.782007324710373
$ws.Range("I24").Value = 3.487691027762498
$ws.Range("J24").Value = 11.00201857869584
$ws.Range("K24").Value = 13.93097215459811
$ws.Range("L24").Value = 10.7992968307174
$ws.Range("M24").Value = 0
$ws.Range("N24").Value = 0
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 0

# Row 25
$ws.Range("B25").Value = 26.52332247880308
$ws.Range("C25").Value = 21.81847615986391
$ws.Range("D25").Value = 10.0767815781245
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 29.45716084570411
$ws.Range("G25").Value = 33.53479946165272
$ws.Range("H25").Value = 2.254401189848858
$ws.Range("I25").Value = 3.084950524584404
$ws.Range("J25").Value = 10.73246904921205
$ws.Range("K25").Value = 14.03569499787289
$ws.Range("L25").Value = 10.06238328661189
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 0

